$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two existing columns (A, B) of data in rows 1-3 with their new values,
# and append two brand-new rows (4 and 5) extending the used range to A1:B5.
$ws.Range("A1").Value = -0.084239552707853341
$ws.Range("B1").Value = -0.084224276384226374
$ws.Range("A2").Value = -0.036072844077985278
$ws.Range("B2").Value = -0.030708633422882319
$ws.Range("A3").Value = -0.013904650964170083
$ws.Range("B3").Value = -0.011197439677351694
$ws.Range("A4").Value = -0.0024887395514233733
$ws.Range("B4").Value = 0.0024887395528284569
$ws.Range("A5").Value = 0.016282035135380166
$ws.Range("B5").Value = -0.016282035208911722

# Column A widened from 14.42578125 to 15.42578125 (+1 char), column B widened
# from 14.42578125 to 14.7109375 (+0.28515625 char).
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334
